# Cross-browser testing support: insert a "browser" column into the DATA
# sheet (between "execute" and "username") and populate it with the
# browsers each test row should run against.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Shift username/password/fristname/lastname columns one to the right to
# make room for the new "browser" column at C.
$ws.Columns.Item(3).Insert()

# Header
$ws.Cells.Item(1, 3).Value = "browser"

# Row 2 (loginLogoutTest / yes) -> chrome
$ws.Cells.Item(2, 3).Value = "chrome"

# Row 3 (loginLogoutTest) -> firefox, and flip its "execute" flag to yes
$ws.Cells.Item(3, 2).Value = "yes"
$ws.Cells.Item(3, 3).Value = "firefox"

# Row 4 (newTest) -> no browser assigned, "execute" flag flips to no
$ws.Cells.Item(4, 2).Value = "no"
$ws.Cells.Item(4, 3).Value = "'"

# Row 5 (loginLogoutTest) -> no browser assigned, "execute" flag flips to no
$ws.Cells.Item(5, 2).Value = "no"
$ws.Cells.Item(5, 3).Value = "'"

# Move the active selection to reflect where the user was working next.
$ws.Range("C6").Select() | Out-Null
